$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing values
$ws.Range("B2").Value = 69
$ws.Range("B3").Value = 56
$ws.Range("B4").Value = 34

# Replace row 5 with what used to be row 7, then delete old rows 5-7 remainder
$ws.Range("A5").Value = "Wrong_Tag_I_as_E"
$ws.Range("B5").Value = 1

# Delete the now-obsolete rows 6 and 7 (shifts cells up, clearing them)
$ws.Range("A6:B7").Delete([Microsoft.Office.Interop.Excel.XlDeleteShiftDirection]::xlShiftUp)
